# Generate Report for handback
# Adds two new rows (6 and 7) to each of the three worksheets
# (Overview, zh-cn, de-de) describing the handback of two new files:
#   2e6d1a54-39a7-4085-a472-e501c0b6163a.md
#   dc1860a5-3157-44b4-96b0-47329f46fa73.md

$wb = $excel.ActiveWorkbook

$file1 = "2e6d1a54-39a7-4085-a472-e501c0b6163a.md"
$file2 = "dc1860a5-3157-44b4-96b0-47329f46fa73.md"

$handed = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = $file1
$wsOverview.Range("B6").Value = $handed
$wsOverview.Range("C6").Value = $handed

$wsOverview.Range("A7").Value = $file2
$wsOverview.Range("B7").Value = $handed
$wsOverview.Range("C7").Value = $handed

$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), $file1, "", "", $file1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), $file2, "", "", $file2)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf1 = "2e6d1a54-39a7-4085-a472-e501c0b6163a.d1a4554e3894a2de5ddd3ad05aef20c9680aa25e.zh-cn.xlf"
$zhXlf2 = "dc1860a5-3157-44b4-96b0-47329f46fa73.4d17e03f44bb8ceffbf727548fcfb0715b0c7c72.zh-cn.xlf"

$wsZh.Range("A6").Value = $file1
$wsZh.Range("B6").Value = $handed
$wsZh.Range("C6").Value = $zhXlf1
$wsZh.Range("D6").Value = "2016-02-15 04:21:01"
$wsZh.Range("E6").Value = $file1
$wsZh.Range("F6").Value = $zhXlf1
$wsZh.Range("G6").Value = "2016-02-15 04:22:06"
$wsZh.Range("H6").Value = "Include"

$wsZh.Range("A7").Value = $file2
$wsZh.Range("B7").Value = $handed
$wsZh.Range("C7").Value = $zhXlf2
$wsZh.Range("D7").Value = "2016-02-15 04:21:01"
$wsZh.Range("E7").Value = $file2
$wsZh.Range("F7").Value = $zhXlf2
$wsZh.Range("G7").Value = "2016-02-15 04:22:06"
$wsZh.Range("H7").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), $file1, "", "", $file1)
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), $zhXlf1, "", "", $zhXlf1)
$wsZh.Hyperlinks.Add($wsZh.Range("E6"), $file1, "", "", $file1)
$wsZh.Hyperlinks.Add($wsZh.Range("F6"), $zhXlf1, "", "", $zhXlf1)

$wsZh.Hyperlinks.Add($wsZh.Range("A7"), $file2, "", "", $file2)
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), $zhXlf2, "", "", $zhXlf2)
$wsZh.Hyperlinks.Add($wsZh.Range("E7"), $file2, "", "", $file2)
$wsZh.Hyperlinks.Add($wsZh.Range("F7"), $zhXlf2, "", "", $zhXlf2)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf1 = "2e6d1a54-39a7-4085-a472-e501c0b6163a.d1a4554e3894a2de5ddd3ad05aef20c9680aa25e.de-de.xlf"
$deXlf2 = "dc1860a5-3157-44b4-96b0-47329f46fa73.4d17e03f44bb8ceffbf727548fcfb0715b0c7c72.de-de.xlf"

$wsDe.Range("A6").Value = $file1
$wsDe.Range("B6").Value = $handed
$wsDe.Range("C6").Value = $deXlf1
$wsDe.Range("D6").Value = "2016-02-15 04:21:16"
$wsDe.Range("E6").Value = $file1
$wsDe.Range("F6").Value = $deXlf1
$wsDe.Range("G6").Value = "2016-02-15 04:22:32"
$wsDe.Range("H6").Value = "Include"

$wsDe.Range("A7").Value = $file2
$wsDe.Range("B7").Value = $handed
$wsDe.Range("C7").Value = $deXlf2
$wsDe.Range("D7").Value = "2016-02-15 04:21:16"
$wsDe.Range("E7").Value = $file2
$wsDe.Range("F7").Value = $deXlf2
$wsDe.Range("G7").Value = "2016-02-15 04:22:32"
$wsDe.Range("H7").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), $file1, "", "", $file1)
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), $deXlf1, "", "", $deXlf1)
$wsDe.Hyperlinks.Add($wsDe.Range("E6"), $file1, "", "", $file1)
$wsDe.Hyperlinks.Add($wsDe.Range("F6"), $deXlf1, "", "", $deXlf1)

$wsDe.Hyperlinks.Add($wsDe.Range("A7"), $file2, "", "", $file2)
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), $deXlf2, "", "", $deXlf2)
$wsDe.Hyperlinks.Add($wsDe.Range("E7"), $file2, "", "", $file2)
$wsDe.Hyperlinks.Add($wsDe.Range("F7"), $deXlf2, "", "", $deXlf2)

Write-Output "Handback report rows added"
